# Updates currentAveragePrice / LevePrice / LeveProfit figures (columns H-N)
# on the Chocobo Gil-profit tracker sheets, per the latest market-board refresh.
# Each block targets one Leve row (matched by its "Leve Item ID" in column G,
# used below only as a human-readable comment) across the eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 399.66666
$ws.Range("I18").Value = 299.5
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 299.5
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -15.5
$ws.Range("N18").Value = -1168

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 3236
$ws.Range("I76").Value = 3239.6
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3239.6
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -2924.6
$ws.Range("N76").Value = -3830

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 3236
$ws.Range("I79").Value = 3239.6
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3239.6
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -2147.6
$ws.Range("N79").Value = -5384

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 1307.2812
$ws.Range("J112").Value = 1361.1
$ws.Range("L112").Value = 4083.3
$ws.Range("N112").Value = -6299.299999999999

# Row 123: Nearly Bare
$ws.Range("H123").Value = 41807.5
$ws.Range("J123").Value = 41807.5
$ws.Range("L123").Value = 41807.5
$ws.Range("N123").Value = -51607.5

# Row 129: Practical Command
$ws.Range("H129").Value = 851.6598
$ws.Range("J129").Value = 903.7727
$ws.Range("L129").Value = 2711.3181
$ws.Range("N129").Value = -12711.3181

# Row 133: Big Brush, Big Dreams
$ws.Range("H133").Value = 45648.89
$ws.Range("J133").Value = 45648.89
$ws.Range("L133").Value = 45648.89
$ws.Range("N133").Value = -55768.89

# Row 138: All-night Crafting
$ws.Range("H138").Value = 6715.5
$ws.Range("I138").Value = 831.1429000000001
$ws.Range("J138").Value = 8279.696
$ws.Range("K138").Value = 2493.4287
$ws.Range("L138").Value = 24839.088
$ws.Range("M138").Value = 2646.5713
$ws.Range("N138").Value = -35119.088

$ws = $wb.Worksheets.Item("ARM")
# Row 3: Skillet Labor
$ws.Range("H3").Value = 13750
$ws.Range("J3").Value = 13750
$ws.Range("L3").Value = 13750
$ws.Range("N3").Value = -13980

# Row 34: Insistent Sallets
$ws.Range("H34").Value = 20471.285
$ws.Range("I34").Value = 38000
$ws.Range("J34").Value = 17549.834
$ws.Range("K34").Value = 38000
$ws.Range("L34").Value = 17549.834
$ws.Range("M34").Value = -37729
$ws.Range("N34").Value = -18091.834

# Row 109: A Head of Demand
$ws.Range("H109").Value = 30900
$ws.Range("J109").Value = 30900
$ws.Range("L109").Value = 30900
$ws.Range("N109").Value = -33674

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1555.2632
$ws.Range("I110").Value = 1469.4166
$ws.Range("J110").Value = 1702.4286
$ws.Range("K110").Value = 1469.4166
$ws.Range("L110").Value = 1702.4286
$ws.Range("M110").Value = 575.5834
$ws.Range("N110").Value = -5792.4286

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2352.9524
$ws.Range("I132").Value = 1377.2354
$ws.Range("J132").Value = 6499.75
$ws.Range("K132").Value = 4131.706200000001
$ws.Range("L132").Value = 19499.25
$ws.Range("M132").Value = -1601.706200000001
$ws.Range("N132").Value = -24559.25

# Row 137: Odd Instruments
$ws.Range("H137").Value = 45146
$ws.Range("J137").Value = 45146
$ws.Range("L137").Value = 45146
$ws.Range("N137").Value = -55346

$ws = $wb.Worksheets.Item("BSM")
# Row 7: Thank You for Your Business
$ws.Range("H7").Value = 19671.75
$ws.Range("I7").Value = 17797.4
$ws.Range("J7").Value = 20523.727
$ws.Range("K7").Value = 17797.4
$ws.Range("L7").Value = 20523.727
$ws.Range("M7").Value = -17684.4
$ws.Range("N7").Value = -20749.727

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 4696417
$ws.Range("I105").Value = 4903540.5
$ws.Range("K105").Value = 4903540.5
$ws.Range("M105").Value = -4901793.5

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2652.5
$ws.Range("I134").Value = 1935.2
$ws.Range("J134").Value = 5214.2856
$ws.Range("K134").Value = 5805.6
$ws.Range("L134").Value = 15642.8568
$ws.Range("M134").Value = -3270.6
$ws.Range("N134").Value = -20712.8568

$ws = $wb.Worksheets.Item("CRP")
# Row 5: Bowing Out
$ws.Range("H5").Value = 3459.4
$ws.Range("I5").Value = 265.66666
$ws.Range("J5").Value = 8250
$ws.Range("K5").Value = 265.66666
$ws.Range("L5").Value = 8250
$ws.Range("M5").Value = -153.66666
$ws.Range("N5").Value = -8474

# Row 17: Say It with Spears
$ws.Range("H17").Value = 15000
$ws.Range("J17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15348

# Row 41: The Lone Bowman
$ws.Range("H41").Value = 31178.857
$ws.Range("I41").Value = 14333.333
$ws.Range("J41").Value = 43813
$ws.Range("K41").Value = 14333.333
$ws.Range("L41").Value = 43813
$ws.Range("M41").Value = -13905.333
$ws.Range("N41").Value = -44669

# Row 50: The Arsenal of Theocracy
$ws.Range("H50").Value = 27993.572
$ws.Range("J50").Value = 27993.572
$ws.Range("L50").Value = 27993.572
$ws.Range("N50").Value = -29243.572

# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 32322.857
$ws.Range("J51").Value = 32322.857
$ws.Range("L51").Value = 32322.857
$ws.Range("N51").Value = -33794.857

# Row 59: Bow Down to Magic
$ws.Range("H59").Value = 28100.8
$ws.Range("I59").Value = 5000
$ws.Range("J59").Value = 30667.555
$ws.Range("K59").Value = 5000
$ws.Range("L59").Value = 30667.555
$ws.Range("M59").Value = -3855
$ws.Range("N59").Value = -32957.555

# Row 60: Bowing to Greater Power
$ws.Range("H60").Value = 20414.58
$ws.Range("J60").Value = 20414.58
$ws.Range("L60").Value = 20414.58
$ws.Range("N60").Value = -21436.58

# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 32322.857
$ws.Range("J61").Value = 32322.857
$ws.Range("L61").Value = 32322.857
$ws.Range("N61").Value = -33018.857

# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 38790.645
$ws.Range("J68").Value = 38790.645
$ws.Range("L68").Value = 38790.645
$ws.Range("N68").Value = -40288.645

# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 38790.645
$ws.Range("J71").Value = 38790.645
$ws.Range("L71").Value = 116371.935
$ws.Range("N71").Value = -123859.935

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1408.381
$ws.Range("I105").Value = 1214.3077
$ws.Range("J105").Value = 1723.75
$ws.Range("K105").Value = 1214.3077
$ws.Range("L105").Value = 1723.75
$ws.Range("M105").Value = 532.6922999999999
$ws.Range("N105").Value = -5217.75

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1904.1052
$ws.Range("I122").Value = 1016
$ws.Range("J122").Value = 3125.25
$ws.Range("K122").Value = 3048
$ws.Range("L122").Value = 9375.75
$ws.Range("M122").Value = -598
$ws.Range("N122").Value = -14275.75

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2544.1724
$ws.Range("I132").Value = 1446.7368
$ws.Range("K132").Value = 4340.2104
$ws.Range("M132").Value = -1810.2104

$ws = $wb.Worksheets.Item("CUL")
# Row 110: His Dark Utensils
$ws.Range("H110").Value = 1527
$ws.Range("I110").Value = 1527
$ws.Range("K110").Value = 4581
$ws.Range("M110").Value = -491

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 4311330.5
$ws.Range("I113").Value = 963.4706
$ws.Range("J113").Value = 10417684
$ws.Range("K113").Value = 2890.4118
$ws.Range("L113").Value = 31253052
$ws.Range("M113").Value = -720.4117999999999
$ws.Range("N113").Value = -31257392

# Row 138: Bring Me Your Tacos
$ws.Range("H138").Value = 3127.7778
$ws.Range("I138").Value = 3037.5
$ws.Range("J138").Value = 3200
$ws.Range("K138").Value = 9112.5
$ws.Range("L138").Value = 9600
$ws.Range("M138").Value = -3972.5
$ws.Range("N138").Value = -19880

# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 2575.6191
$ws.Range("I140").Value = 2875.7646
$ws.Range("J140").Value = 1300
$ws.Range("K140").Value = 8627.293799999999
$ws.Range("L140").Value = 3900
$ws.Range("M140").Value = -3447.293799999999
$ws.Range("N140").Value = -14260

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me
$ws.Range("H5").Value = 12994.25
$ws.Range("J5").Value = 12994.25
$ws.Range("L5").Value = 12994.25
$ws.Range("N5").Value = -13218.25

# Row 14: All That Glitters
$ws.Range("H14").Value = 16500075
$ws.Range("I14").Value = 16500075
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 16500075
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -16499907
$ws.Range("N14").ClearContents()

# Row 22: Bad to the Bone
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 31744.857
$ws.Range("J46").Value = 31744.857
$ws.Range("L46").Value = 31744.857
$ws.Range("N46").Value = -32056.857

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 35716744
$ws.Range("I80").Value = 50002340
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 50002340
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -50001342
$ws.Range("N80").Value = -4746

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 35716744
$ws.Range("I83").Value = 50002340
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 250011700
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -250006708
$ws.Range("N83").Value = -23734

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 7500.5713
$ws.Range("I122").Value = 3165.3333
$ws.Range("J122").Value = 10752
$ws.Range("K122").Value = 9495.999899999999
$ws.Range("L122").Value = 32256
$ws.Range("M122").Value = -7045.999899999999
$ws.Range("N122").Value = -37156

# Row 139: Ringing Gratitude
$ws.Range("H139").Value = 37857.145
$ws.Range("J139").Value = 37857.145
$ws.Range("L139").Value = 37857.145
$ws.Range("N139").Value = -48137.145

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 1343.3334
$ws.Range("I16").Value = 1343.3334
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1343.3334
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1173.3334
$ws.Range("N16").ClearContents()

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 655.1111
$ws.Range("I68").Value = 655.1111
$ws.Range("K68").Value = 655.1111
$ws.Range("M68").Value = 93.88890000000004

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 655.1111
$ws.Range("I71").Value = 655.1111
$ws.Range("K71").Value = 3275.5555
$ws.Range("M71").Value = 468.4445000000001

# Row 122: Hell on Leather
$ws.Range("H122").Value = 6253.857
$ws.Range("J122").Value = 7705
$ws.Range("L122").Value = 23115
$ws.Range("N122").Value = -28015

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 252625250
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
